$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "v1.2.2-0 (c17b9bc9)" run results: fill in columns G (Seconds) / H (Ratio) ---
$ws.Range("G7").Value  = 12673
$ws.Range("H7").Value  = 0.14000000000000001

$ws.Range("G8").Value  = 10375
$ws.Range("H8").Value  = 0.38400000000000001

$ws.Range("G9").Value  = 4774
$ws.Range("H9").Value  = 0.98899999999999999

$ws.Range("G10").Value = 10257
$ws.Range("H10").Value = 0.159

$ws.Range("G11").Value = 5285
$ws.Range("H11").Value = 0.26800000000000002

$ws.Range("G12").Value = 4305
$ws.Range("H12").Value = 0.41

$ws.Range("G13").Value = 10327
$ws.Range("H13").Value = 0.13900000000000001

$ws.Range("G14").Value = 6737
$ws.Range("H14").Value = 0.253

$ws.Range("G15").Value = 4807
$ws.Range("H15").Value = 0.61199999999999999

$ws.Range("G16").Value = 13272
$ws.Range("H16").Value = 0.46600000000000003

# --- Base Score row: new run's base score + DBT:QEMU ratio formula ---
$ws.Range("G18").Value = 0.31390000000000001
$ws.Range("I18").Formula = '=$B$18/G18'

# --- Format the Factor column for the Base Score row with 4 decimals ---
$ws.Range("F18").NumberFormat = "0.0000"
$ws.Range("I18").NumberFormat = "0.0000"

# --- Selection state ---
$excel.Goto($ws.Range("G18:H18"))
